# Update the crypto symbol list (GitHub Actions scraper refresh).
# Columns D (Price) and E (Volume(1h)) hold numeric-looking text
# (e.g. "331.81", "0.90%") that must stay plain text, exactly as the
# original workbook stored them (inline strings, not numbers/percent
# values). Prefixing the literal with a single quote makes Excel apply
# its normal "quote prefix" text coercion so the value is kept verbatim
# instead of being parsed into a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''331.81'
$ws.Range("E2").Value = '''0.90%'
$ws.Range("D3").Value = '''44.69'
$ws.Range("E3").Value = '''1.07%'
$ws.Range("D4").Value = '''5.551'
$ws.Range("E4").Value = '''-0.49%'
$ws.Range("D5").Value = '''0.08187'
$ws.Range("E5").Value = '''1.44%'
$ws.Range("D6").Value = '''2.058'
$ws.Range("E6").Value = '''4.35%'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '''0.9761'
$ws.Range("E7").Value = '''2.49%'
$ws.Range("B8").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C8").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D8").Value = '''0.1117'
$ws.Range("E8").Value = '''-3.51%'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '''0.1901'
$ws.Range("E9").Value = '''2.53%'
$ws.Range("B10").Value = 'MCDex'
$ws.Range("C10").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D10").Value = '''10.20'
$ws.Range("E10").Value = '''-13.95%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.1005'
$ws.Range("E11").Value = '''2.65%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.04722'
$ws.Range("E12").Value = '''-0.22%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.1058'
$ws.Range("E13").Value = '''-0.98%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001263'
$ws.Range("E14").Value = '''-1.77%'
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").Value = '''0.04110'
$ws.Range("E15").Value = '''-2.94%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005961'
$ws.Range("E16").Value = '''-0.21%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.350'
$ws.Range("E17").Value = '''-0.64%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''4.430'
$ws.Range("E18").Value = '''2.38%'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '''2.645'
$ws.Range("E19").Value = '''2.89%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3351'
$ws.Range("E20").Value = '''-3.53%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '''0.1392'
$ws.Range("E21").Value = '''-1.19%'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").Value = '''0.2566'
$ws.Range("E22").Value = '''2.35%'
$ws.Range("B23").Value = 'BitKan'
$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D23").Value = '''0.001301'
$ws.Range("E23").Value = '''3.76%'
$ws.Range("B24").Value = 'HotbitToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D24").Value = '''0.004418'
$ws.Range("E24").Value = '''2.27%'
$ws.Range("E38").Value = '''3.51%'
$ws.Range("D39").Value = '''0.05728'
$ws.Range("E39").Value = '''3.41%'
$ws.Range("D40").Value = '''0.007625'
$ws.Range("E40").Value = '''0.93%'
$ws.Range("D41").Value = '''0.1421'
$ws.Range("E41").Value = '''0.90%'
$ws.Range("D42").Value = '''0.007551'
$ws.Range("E42").Value = '''-6.92%'
$ws.Range("E43").Value = '''-3.06%'
$ws.Range("E44").Value = '''-6.59%'
$ws.Range("D45").Value = '''0.00007049'
$ws.Range("E45").Value = '''-1.45%'
$ws.Range("E46").Value = '''-0.25%'
$ws.Range("D47").Value = '''0.0005793'
$ws.Range("E47").Value = '''-0.32%'
$ws.Range("D48").Value = '''0.002516'
$ws.Range("E48").Value = '''9.30%'
$ws.Range("D49").Value = '''0.003537'
$ws.Range("E49").Value = '''0.31%'
$ws.Range("D50").Value = '''0.00002097'
$ws.Range("E50").Value = '''-0.25%'
$ws.Range("D51").Value = '''0.0001997'
$ws.Range("E51").Value = '''-0.25%'
